# Added: "list" sheet (per-doctor patient list headers), kept "form" first.
# Removed tabSelected from "form" (now "list" is the active tab).
# Also nudge the "form" view scroll (topLeftCell A5 -> A7, best effort).

$wb = $excel.ActiveWorkbook

# Add the new "list" worksheet and place it right after "form"
# (sheet references can be reseated by Add()/Move(), so always resolve
# the target sheet by name immediately before using it)
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "list"
$newSheet.Move($null, $wb.Worksheets.Item("form"))

# best-effort scroll position nudge on the form sheet (A5 -> A7)
$formSheet = $wb.Worksheets.Item("form")
$formSheet.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1

# Header row 1: doc id / doc_name
$listSheet = $wb.Worksheets.Item("list")
$listSheet.Range("A1").Value = "doc id="
$listSheet.Range("C1").Value = "doc_name"

# Header row 2: patient num / patient name / sex / birth date
$listSheet.Range("A2").Value = "patient num"
$listSheet.Range("B2").Value = "patient name"
$listSheet.Range("C2").Value = "М\Ж\Р"
$listSheet.Range("D2").Value = "Дата рождения"

# column widths (best effort; engine quantizes to 1/6 character units)
$listSheet.Columns.Item(1).ColumnWidth = 10.166666666666666
$listSheet.Columns.Item(2).ColumnWidth = 11.072916666666666
$listSheet.Columns.Item(3).ColumnWidth = 8.709635416666666
$listSheet.Columns.Item(4).ColumnWidth = 13.346354166666666

# make "list" the active/visible tab, with G7 selected
$listSheet.Activate()
$listSheet.Range("G7").Select() | Out-Null
